# RGB Werte.xlsx - add a new "Schlüssel" (key) entry as row 15
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A15").Value = "Schlüssel"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1

# Restore the previously-selected cell/view
$ws.Range("D15").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
